$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45194
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/caja 14 kilos granel'
$ws.Range("S2").Value = 1071

$ws.Range("D3").Value = 45238
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 14 kilos empedrada'
$ws.Range("S3").Value = 1286

$ws.Range("D4").Value = 45236
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/caja 14 kilos granel'
$ws.Range("S4").Value = 1286

$ws.Range("D6").Value = 45212
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = '$/caja 14 kilos granel'
$ws.Range("S6").Value = 1214

$ws.Range("D7").Value = 45240
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = '$/caja 14 kilos granel'
$ws.Range("S7").Value = 1143

$ws.Range("D8").Value = 45196
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 14 kilos empedrada'
$ws.Range("S8").Value = 1071

$ws.Range("D9").Value = 45167
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 16000
$ws.Range("Q9").Value = '$/caja 14 kilos empedrada'
$ws.Range("S9").Value = 1143

$ws.Range("D10").Value = 44210
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 10357
$ws.Range("Q10").Value = '$/caja 14 kilos empedrada'
$ws.Range("S10").Value = 740

$ws.Range("D11").Value = 44229
$ws.Range("M11").Value = 55
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11364
$ws.Range("Q11").Value = '$/caja 14 kilos empedrada'
$ws.Range("S11").Value = 812

$ws.Range("D12").Value = 45152
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("Q12").Value = '$/caja 14 kilos empedrada'
$ws.Range("S12").Value = 1143

$ws.Range("D13").Value = 44172
$ws.Range("M13").Value = 90
$ws.Range("N13").Value = 8500
$ws.Range("O13").Value = 9000
$ws.Range("P13").Value = 8806
$ws.Range("Q13").Value = '$/caja 14 kilos empedrada'
$ws.Range("S13").Value = 629

$ws.Range("D14").Value = 45224
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = '$/caja 14 kilos granel'
$ws.Range("S14").Value = 1071

$ws.Range("D15").Value = 45138
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = '$/caja 14 kilos granel'
$ws.Range("S15").Value = 1000

$ws.Range("D16").Value = 44253
$ws.Range("M16").Value = 90
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 13000
$ws.Range("P16").Value = 12667
$ws.Range("Q16").Value = '$/caja 14 kilos empedrada'
$ws.Range("S16").Value = 905

$ws.Range("D17").Value = 44216
$ws.Range("M17").Value = 55
$ws.Range("N17").Value = 11000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 11545
$ws.Range("Q17").Value = '$/caja 14 kilos empedrada'
$ws.Range("S17").Value = 825

$ws.Range("D20").Value = 45155
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 15000
$ws.Range("Q20").Value = '$/caja 14 kilos empedrada'
$ws.Range("S20").Value = 1071

$ws.Range("D21").Value = 45140
$ws.Range("M21").Value = 30
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("Q21").Value = '$/caja 14 kilos granel'
$ws.Range("S21").Value = 1071

$ws.Range("D22").Value = 44232
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 11000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 11583
$ws.Range("Q22").Value = '$/caja 14 kilos empedrada'
$ws.Range("S22").Value = 827

$ws.Range("D23").Value = 44181
$ws.Range("M23").Value = 65
$ws.Range("N23").Value = 9000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 9462
$ws.Range("Q23").Value = '$/caja 14 kilos empedrada'
$ws.Range("S23").Value = 676

$ws.Range("D24").Value = 45222
$ws.Range("M24").Value = 80
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("Q24").Value = '$/caja 14 kilos granel'
$ws.Range("S24").Value = 1071
